$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N builds a semicolon-joined "CSV row" string out of columns A:K for
# each record (node;commodity;level;year;time;value;unit;share;SSP5 value;Unit).
# This is the helper column used to later paste-as-values into .csv files for
# scenario development (see commit message). Row 14 mirrors the header row
# (row 1); rows 15-56 mirror the data rows 2-43 respectively.
$ws.Range("N14").Formula = '=A1&";"&B1&";"&C1&";"&D1&";"&E1&";"&F1&";"&G1&";"&H1&";"&I1&";"&J1&";"&K1'
$ws.Range("N15").Formula = '=A2&";"&B2&";"&C2&";"&D2&";"&E2&";"&F2&";"&G2&";"&H2&";"&I2&";"&J2&";"&K2'
$ws.Range("N16").Formula = '=A3&";"&B3&";"&C3&";"&D3&";"&E3&";"&F3&";"&G3&";"&H3&";"&I3&";"&J3&";"&K3'
$ws.Range("N17").Formula = '=A4&";"&B4&";"&C4&";"&D4&";"&E4&";"&F4&";"&G4&";"&H4&";"&I4&";"&J4&";"&K4'
$ws.Range("N18").Formula = '=A5&";"&B5&";"&C5&";"&D5&";"&E5&";"&F5&";"&G5&";"&H5&";"&I5&";"&J5&";"&K5'
$ws.Range("N19").Formula = '=A6&";"&B6&";"&C6&";"&D6&";"&E6&";"&F6&";"&G6&";"&H6&";"&I6&";"&J6&";"&K6'
$ws.Range("N20").Formula = '=A7&";"&B7&";"&C7&";"&D7&";"&E7&";"&F7&";"&G7&";"&H7&";"&I7&";"&J7&";"&K7'
$ws.Range("N21").Formula = '=A8&";"&B8&";"&C8&";"&D8&";"&E8&";"&F8&";"&G8&";"&H8&";"&I8&";"&J8&";"&K8'
$ws.Range("N22").Formula = '=A9&";"&B9&";"&C9&";"&D9&";"&E9&";"&F9&";"&G9&";"&H9&";"&I9&";"&J9&";"&K9'
$ws.Range("N23").Formula = '=A10&";"&B10&";"&C10&";"&D10&";"&E10&";"&F10&";"&G10&";"&H10&";"&I10&";"&J10&";"&K10'
$ws.Range("N24").Formula = '=A11&";"&B11&";"&C11&";"&D11&";"&E11&";"&F11&";"&G11&";"&H11&";"&I11&";"&J11&";"&K11'
$ws.Range("N25").Formula = '=A12&";"&B12&";"&C12&";"&D12&";"&E12&";"&F12&";"&G12&";"&H12&";"&I12&";"&J12&";"&K12'
$ws.Range("N26").Formula = '=A13&";"&B13&";"&C13&";"&D13&";"&E13&";"&F13&";"&G13&";"&H13&";"&I13&";"&J13&";"&K13'
$ws.Range("N27").Formula = '=A14&";"&B14&";"&C14&";"&D14&";"&E14&";"&F14&";"&G14&";"&H14&";"&I14&";"&J14&";"&K14'
$ws.Range("N28").Formula = '=A15&";"&B15&";"&C15&";"&D15&";"&E15&";"&F15&";"&G15&";"&H15&";"&I15&";"&J15&";"&K15'
$ws.Range("N29").Formula = '=A16&";"&B16&";"&C16&";"&D16&";"&E16&";"&F16&";"&G16&";"&H16&";"&I16&";"&J16&";"&K16'
$ws.Range("N30").Formula = '=A17&";"&B17&";"&C17&";"&D17&";"&E17&";"&F17&";"&G17&";"&H17&";"&I17&";"&J17&";"&K17'
$ws.Range("N31").Formula = '=A18&";"&B18&";"&C18&";"&D18&";"&E18&";"&F18&";"&G18&";"&H18&";"&I18&";"&J18&";"&K18'
$ws.Range("N32").Formula = '=A19&";"&B19&";"&C19&";"&D19&";"&E19&";"&F19&";"&G19&";"&H19&";"&I19&";"&J19&";"&K19'
$ws.Range("N33").Formula = '=A20&";"&B20&";"&C20&";"&D20&";"&E20&";"&F20&";"&G20&";"&H20&";"&I20&";"&J20&";"&K20'
$ws.Range("N34").Formula = '=A21&";"&B21&";"&C21&";"&D21&";"&E21&";"&F21&";"&G21&";"&H21&";"&I21&";"&J21&";"&K21'
$ws.Range("N35").Formula = '=A22&";"&B22&";"&C22&";"&D22&";"&E22&";"&F22&";"&G22&";"&H22&";"&I22&";"&J22&";"&K22'
$ws.Range("N36").Formula = '=A23&";"&B23&";"&C23&";"&D23&";"&E23&";"&F23&";"&G23&";"&H23&";"&I23&";"&J23&";"&K23'
$ws.Range("N37").Formula = '=A24&";"&B24&";"&C24&";"&D24&";"&E24&";"&F24&";"&G24&";"&H24&";"&I24&";"&J24&";"&K24'
$ws.Range("N38").Formula = '=A25&";"&B25&";"&C25&";"&D25&";"&E25&";"&F25&";"&G25&";"&H25&";"&I25&";"&J25&";"&K25'
$ws.Range("N39").Formula = '=A26&";"&B26&";"&C26&";"&D26&";"&E26&";"&F26&";"&G26&";"&H26&";"&I26&";"&J26&";"&K26'
$ws.Range("N40").Formula = '=A27&";"&B27&";"&C27&";"&D27&";"&E27&";"&F27&";"&G27&";"&H27&";"&I27&";"&J27&";"&K27'
$ws.Range("N41").Formula = '=A28&";"&B28&";"&C28&";"&D28&";"&E28&";"&F28&";"&G28&";"&H28&";"&I28&";"&J28&";"&K28'
$ws.Range("N42").Formula = '=A29&";"&B29&";"&C29&";"&D29&";"&E29&";"&F29&";"&G29&";"&H29&";"&I29&";"&J29&";"&K29'
$ws.Range("N43").Formula = '=A30&";"&B30&";"&C30&";"&D30&";"&E30&";"&F30&";"&G30&";"&H30&";"&I30&";"&J30&";"&K30'
$ws.Range("N44").Formula = '=A31&";"&B31&";"&C31&";"&D31&";"&E31&";"&F31&";"&G31&";"&H31&";"&I31&";"&J31&";"&K31'
$ws.Range("N45").Formula = '=A32&";"&B32&";"&C32&";"&D32&";"&E32&";"&F32&";"&G32&";"&H32&";"&I32&";"&J32&";"&K32'
$ws.Range("N46").Formula = '=A33&";"&B33&";"&C33&";"&D33&";"&E33&";"&F33&";"&G33&";"&H33&";"&I33&";"&J33&";"&K33'
$ws.Range("N47").Formula = '=A34&";"&B34&";"&C34&";"&D34&";"&E34&";"&F34&";"&G34&";"&H34&";"&I34&";"&J34&";"&K34'
$ws.Range("N48").Formula = '=A35&";"&B35&";"&C35&";"&D35&";"&E35&";"&F35&";"&G35&";"&H35&";"&I35&";"&J35&";"&K35'
$ws.Range("N49").Formula = '=A36&";"&B36&";"&C36&";"&D36&";"&E36&";"&F36&";"&G36&";"&H36&";"&I36&";"&J36&";"&K36'
$ws.Range("N50").Formula = '=A37&";"&B37&";"&C37&";"&D37&";"&E37&";"&F37&";"&G37&";"&H37&";"&I37&";"&J37&";"&K37'
$ws.Range("N51").Formula = '=A38&";"&B38&";"&C38&";"&D38&";"&E38&";"&F38&";"&G38&";"&H38&";"&I38&";"&J38&";"&K38'
$ws.Range("N52").Formula = '=A39&";"&B39&";"&C39&";"&D39&";"&E39&";"&F39&";"&G39&";"&H39&";"&I39&";"&J39&";"&K39'
$ws.Range("N53").Formula = '=A40&";"&B40&";"&C40&";"&D40&";"&E40&";"&F40&";"&G40&";"&H40&";"&I40&";"&J40&";"&K40'
$ws.Range("N54").Formula = '=A41&";"&B41&";"&C41&";"&D41&";"&E41&";"&F41&";"&G41&";"&H41&";"&I41&";"&J41&";"&K41'
$ws.Range("N55").Formula = '=A42&";"&B42&";"&C42&";"&D42&";"&E42&";"&F42&";"&G42&";"&H42&";"&I42&";"&J42&";"&K42'
$ws.Range("N56").Formula = '=A43&";"&B43&";"&C43&";"&D43&";"&E43&";"&F43&";"&G43&";"&H43&";"&I43&";"&J43&";"&K43'

# Reflect the author's final selection (the newly filled N14:N56 helper column).
$ws.Range("N14:N56").Select()
